# "Blog con APP_revisando AVATAR"
# Replace the old "Casos de prueba" test-case table (rows 7-14, with a
# duplicated header at row 11) with a new, single, 12-row test case table
# (rows 7-18), drop the stray formatting that used to spill into columns
# F:X, widen column B a bit, and refresh the sheet view / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Wipe the old test-case area completely (values + formatting), from
#    row 7 all the way down to the bottom of the old used range. This
#    removes: the 3 original case rows (7-9), the stray per-cell styling
#    that used to extend out to column X, the duplicated header row
#    (old row 11) and the 3 rows of data that used to follow it
#    (old rows 12-14).
# ---------------------------------------------------------------------
$ws.Range("A7:X1000").Clear()

# ---------------------------------------------------------------------
# 2. Write the refreshed set of test cases (Caso #1 .. Caso #12).
# ---------------------------------------------------------------------
$colA = @("Caso #1","Caso #2","Caso #3","Caso #4","Caso #5","Caso #6","Caso #7","Caso #8","Caso #9","Caso #10","Caso #11","Caso #12")
$colB = @(
  "Instalación de Boostrap, descargamos el Componente NAVBar, Card",
  "En posts_page: Error en {% for post in post %}",
  "Problema al no cargar una imagen de portada",
  "En TemplateDoesNotExist at /buscarBlog/ busquedaBlog.html",
  "resultadoBusquedaBlog.html --> arroja la búsqueda como TABLA",
  "En home page agregue un botón de acceso rápido a las búsquedas de Blogs",
  "Problemas al crear la nueva APP Account para Register, Login, Logout, Editar Perfil",
  "imagen_portada = models.ImageField(null=True, blank=True, default = 'default-image.png')",
  "Se creo el servicio de Mensajería sobre los Blogs",
  "Problemas con GITHUB",
  "Problemas con el texto que se muestra en el contenido del Blog.",
  "Problemas con el texto del servicio de Mensajería. No se muestran los comnetarios en los Blogs."
)
$colC = @(
  "No toma el boostrap utilizado",
  "Con ayuda de Coder Ask solucionamos el Error en {% for post in post %}",
  'Se soluciono poniendo por default = "default-image.png"',
  "Error al realizar la búsqueda de un Blog ---> dos errores en la línea 93 de views.py ---> Se Soluciono con: nombre=Post.objects.filter(title__icontains = nombre)",
  "Se soluciono aplicando boostrap,css,js",
  "Error al realizar la búsqueda de un Blog --> colucionado redireccionando a resultadoBusquedaBlog.html",
  "Revisamos con CoderAsk el día completo con 3 Tutores y no encuentran la falla",
  "Había una imagen por defecto en el caso que el usuario no cargue una, por algún motivo ahora no aparece",
  "Me aparece el botón Comentar, sin embargo no me aparece el campo para dejar comentario",
  "Estoy teniendo probelmas para subir al repo los commit, lo estamos revisndo con CoderAsk desde esta mañana. SIN SOLUCIÓN",
  'Se soluciono en "post.html" --> poniendo <p> {{post.description|safe}} </p>',
  $null
)
$colD = @(45041,45048,45048,45049,45049,45049,45049,45049,45050,45050,45051,45051)
$colE = @("SI","SI","SI","SI","SI","SI","SI","NO","SI","SI","SI","SI")
$rowHeights = @(30,30,30,60,30,45,30,45,30,45,30,45)

for ($i = 0; $i -lt 12; $i++) {
  $r = 7 + $i
  $ws.Cells.Item($r,1).Value2 = $colA[$i]
  $ws.Cells.Item($r,2).Value2 = $colB[$i]
  if ($colC[$i] -ne $null) {
    $ws.Cells.Item($r,3).Value2 = $colC[$i]
  }
  $ws.Cells.Item($r,4).Value2 = $colD[$i]
  $ws.Cells.Item($r,4).NumberFormat = "d-mmm"
  $ws.Cells.Item($r,5).Value2 = $colE[$i]

  $rowRange = $ws.Range("A" + $r + ":E" + $r)
  $rowRange.Borders.LineStyle = 1
  $rowRange.Borders.Color = 0
  $rowRange.Font.Name = "Calibri"
  $rowRange.Font.Size = 11
  $rowRange.HorizontalAlignment = -4131
  $rowRange.VerticalAlignment = -4130

  $ws.Range("B" + $r + ":C" + $r).WrapText = $true
  $ws.Cells.Item($r,4).WrapText = $true

  $ws.Rows.Item($r).RowHeight = $rowHeights[$i]
}

# ---------------------------------------------------------------------
# 3. Row 19 is blank underneath the table (leftover formatting only on
#    B19), then 20+ carry on as plain 15.75pt filler rows, same as the
#    rest of the sheet.
# ---------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 15.75
$ws.Rows.Item(20).RowHeight = 15.75

# ---------------------------------------------------------------------
# 4. Drop the 4 trailing filler rows (997-1000) that used to pad the
#    sheet out to row 1000; the refreshed sheet only pads out to 996.
# ---------------------------------------------------------------------
$ws.Range("A997:A1000").EntireRow.Delete()

# ---------------------------------------------------------------------
# 5. Column B grew a bit wider.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 44

# ---------------------------------------------------------------------
# 6. Refresh the view: scrolled down a touch, new selection.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E19").Select()
